# Fix #517: a table cell that starts with an (intended) empty leading
# paragraph ends up with a genuinely empty paragraph at the top of the
# cell. Remove that leading empty paragraph from every cell of the
# table so the cell's real content paragraph becomes the first one.
$d = $word.ActiveDocument

foreach ($t in $d.Tables) {
    for ($r = 1; $r -le $t.Rows.Count; $r++) {
        for ($c = 1; $c -le $t.Columns.Count; $c++) {
            $cell = $t.Cell($r, $c)
            # Only strip the leading paragraph when it is empty and the
            # cell has more than one paragraph (so we never eat the only
            # paragraph left in the cell).
            if ($cell.Range.Paragraphs.Count -gt 1) {
                $firstPara = $cell.Range.Paragraphs(1)
                if ($firstPara.Range.Text -eq "`r") {
                    $firstPara.Range.Delete()
                }
            }
        }
    }
}
